$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp header ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Julio de 2020 a las 12:06"

# --- Update per-country statistics with freshly scraped figures ---
$ws.Range("B4").Value = 3414042
$ws.Range("C4").Value = 47
$ws.Range("D4").Value = 1517427
$ws.Range("E4").Value = 1758831
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 137784
$ws.Range("B6").Value = 879888
$ws.Range("C6").Value = 422
$ws.Range("D6").Value = 554888
$ws.Range("E6").Value = 301800
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 13
$ws.Range("H6").Value = 23200
$ws.Range("B19").Value = 199968
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 184600
$ws.Range("E19").Value = 6233
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 9135
$ws.Range("B20").Value = 186894
$ws.Range("C20").Value = 3099
$ws.Range("D20").Value = 98317
$ws.Range("E20").Value = 86186
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 39
$ws.Range("H20").Value = 2391
$ws.Range("B29").Value = 76981
$ws.Range("C29").Value = 1282
$ws.Range("D29").Value = 36689
$ws.Range("E29").Value = 36636
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 50
$ws.Range("H29").Value = 3656
$ws.Range("B35").Value = 57006
$ws.Range("C35").Value = 2784
$ws.Range("D35").Value = 20371
$ws.Range("E35").Value = 35036
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 227
$ws.Range("H35").Value = 1599
$ws.Range("B36").Value = 58179
$ws.Range("C36").Value = 2164
$ws.Range("D36").Value = 37257
$ws.Range("E36").Value = 20663
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 2
$ws.Range("H36").Value = 259
$ws.Range("B37").Value = 55508
$ws.Range("C37").Value = 614
$ws.Range("D37").Value = 45356
$ws.Range("E37").Value = 9759
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = 393
$ws.Range("B43").Value = 46283
$ws.Range("C43").Value = 322
$ws.Range("D43").Value = 42285
$ws.Range("E43").Value = 3972
$ws.Range("F43").Value = 0
$ws.Range("B52").Value = 32948
$ws.Range("C52").Value = 413
$ws.Range("D52").Value = 21692
$ws.Range("E52").Value = 9355
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 17
$ws.Range("H52").Value = 1901
$ws.Range("B65").Value = 15821
$ws.Range("C65").Value = 76
$ws.Range("D65").Value = 12676
$ws.Range("E65").Value = 2892
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 3
$ws.Range("H65").Value = 253
$ws.Range("B76").Value = 9978
$ws.Range("C76").Value = 304
$ws.Range("D76").Value = 5732
$ws.Range("E76").Value = 3979
$ws.Range("F76").Value = 0
$ws.Range("B79").Value = 8725
$ws.Range("C79").Value = 7
$ws.Range("D79").Value = 8520
$ws.Range("E79").Value = 83
$ws.Range("F79").Value = 0
$ws.Range("B85").Value = 7295
$ws.Range("C85").Value = 1
$ws.Range("D85").Value = 6800
$ws.Range("E85").Value = 166
$ws.Range("F85").Value = 0
$ws.Range("B90").Value = 6473
$ws.Range("C90").Value = 243
$ws.Range("D90").Value = 1084
$ws.Range("E90").Value = 5351
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 2
$ws.Range("H90").Value = 38
$ws.Range("B102").Value = 3571
$ws.Range("C102").Value = 117
$ws.Range("D102").Value = 2014
$ws.Range("E102").Value = 1462
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 2
$ws.Range("H102").Value = 95
$ws.Range("B123").Value = 1849
$ws.Range("C123").Value = 8
$ws.Range("D123").Value = 1473
$ws.Range("E123").Value = 265
$ws.Range("F123").Value = 0
$ws.Range("B127").Value = 1522
$ws.Range("C127").Value = 52
$ws.Range("D127").Value = 1217
$ws.Range("E127").Value = 298
$ws.Range("F127").Value = 0


# --- Re-sort the country table (rows 4-219) by "Casos totales" (column B)
#     descending, breaking ties alphabetically by country name (column A) ---
$dataRange = $ws.Range("A4:H219")
$key1 = $ws.Range("B4:B219")
$key2 = $ws.Range("A4:A219")
$dataRange.Sort($key1, 2, $key2, $null, 1, $null, $null, 1)
